# Fix header template import
# The "PO - TS" sheet (3rd tab) has a set of column headers under the
# "Chi phí / ..." (Cost) section that were actually meant to read
# "Chi phí ước tính / ..." (Estimated cost). Also the "Số lượng quy đổi"
# header in that same row was mislabeled and should read "Tỷ lệ quy đổi".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO - TS")

$ws.Range("S1").Value  = "Chi tiết đơn hàng / Tỷ lệ quy đổi (*)"
$ws.Range("Z1").Value  = "Chi phí ước tính / Mã Sản phẩm"
$ws.Range("AA1").Value = "Chi phí ước tính / Tiền tệ"
$ws.Range("AB1").Value = "Chi phí ước tính / Tỷ giá"
$ws.Range("AC1").Value = "Chi phí ước tính / Tổng tiền ngoại tệ̣"
$ws.Range("AD1").Value = "Chi phí ước tính / Thành tiền VND"
$ws.Range("AE1").Value = "Chi phí ước tính / Chi phí trước thuế"

# Restore the view/selection state recorded for this sheet after the edit.
$ws.Activate()
$ws.Range("Z7").Select()
